$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 3929
$ws.Cells.Item(76, 9).Value = 2763.3333
$ws.Cells.Item(76, 10).Value = 4428.5713
$ws.Cells.Item(76, 11).Value = 2763.3333
$ws.Cells.Item(76, 12).Value = 4428.5713
$ws.Cells.Item(76, 13).Value = -2448.3333
$ws.Cells.Item(76, 14).Value = -5058.5713

$ws.Cells.Item(79, 8).Value = 3929
$ws.Cells.Item(79, 9).Value = 2763.3333
$ws.Cells.Item(79, 10).Value = 4428.5713
$ws.Cells.Item(79, 11).Value = 2763.3333
$ws.Cells.Item(79, 12).Value = 4428.5713
$ws.Cells.Item(79, 13).Value = -1671.3333
$ws.Cells.Item(79, 14).Value = -6612.5713

$ws.Cells.Item(100, 8).Value = 5575.5713
$ws.Cells.Item(100, 9).Value = 1428.1818
$ws.Cells.Item(100, 10).Value = 8259.177
$ws.Cells.Item(100, 11).Value = 1428.1818
$ws.Cells.Item(100, 12).Value = 8259.177
$ws.Cells.Item(100, 13).Value = -887.1818000000001
$ws.Cells.Item(100, 14).Value = -9341.177

$ws.Cells.Item(113, 8).Value = 2856.5293
$ws.Cells.Item(113, 9).Value = 1890
$ws.Cells.Item(113, 10).Value = 2985.4
$ws.Cells.Item(113, 11).Value = 1890
$ws.Cells.Item(113, 12).Value = 2985.4
$ws.Cells.Item(113, 13).Value = 1364
$ws.Cells.Item(113, 14).Value = -9493.4

$ws.Cells.Item(125, 8).Value = 3100
$ws.Cells.Item(125, 9).Value = 0
$ws.Cells.Item(125, 10).Value = 3100
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(125, 12).Value = 27900
$ws.Cells.Item(125, 14).Value = -32820
$ws.Cells.Item(125, 13).ClearContents()

$ws.Cells.Item(135, 8).Value = 705.5
$ws.Cells.Item(135, 9).Value = 705.5
$ws.Cells.Item(135, 11).Value = 6349.5
$ws.Cells.Item(135, 13).Value = -3814.5

$ws.Cells.Item(137, 8).Value = 1673
$ws.Cells.Item(137, 9).Value = 1382.125
$ws.Cells.Item(137, 10).Value = 4000
$ws.Cells.Item(137, 11).Value = 4146.375
$ws.Cells.Item(137, 12).Value = 12000
$ws.Cells.Item(137, 13).Value = -1596.375
$ws.Cells.Item(137, 14).Value = -17100

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 963.6667
$ws.Cells.Item(4, 9).Value = 995.5
$ws.Cells.Item(4, 10).Value = 900
$ws.Cells.Item(4, 11).Value = 995.5
$ws.Cells.Item(4, 12).Value = 900
$ws.Cells.Item(4, 13).Value = -879.5
$ws.Cells.Item(4, 14).Value = -1132

$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 14).ClearContents()

$ws.Cells.Item(45, 8).Value = 1389.1875
$ws.Cells.Item(45, 9).Value = 1335.2858
$ws.Cells.Item(45, 10).Value = 1431.1111
$ws.Cells.Item(45, 11).Value = 1335.2858
$ws.Cells.Item(45, 12).Value = 1431.1111
$ws.Cells.Item(45, 13).Value = -958.2858000000001
$ws.Cells.Item(45, 14).Value = -2185.1111

$ws.Cells.Item(74, 8).Value = 12502781
$ws.Cells.Item(74, 9).Value = 19567934
$ws.Cells.Item(74, 10).Value = 2893.2307
$ws.Cells.Item(74, 11).Value = 19567934
$ws.Cells.Item(74, 12).Value = 2893.2307
$ws.Cells.Item(74, 13).Value = -19567060
$ws.Cells.Item(74, 14).Value = -4641.2307

$ws.Cells.Item(77, 8).Value = 12502781
$ws.Cells.Item(77, 9).Value = 19567934
$ws.Cells.Item(77, 10).Value = 2893.2307
$ws.Cells.Item(77, 11).Value = 97839670
$ws.Cells.Item(77, 12).Value = 14466.1535
$ws.Cells.Item(77, 13).Value = -97835302
$ws.Cells.Item(77, 14).Value = -23202.1535

$ws.Cells.Item(102, 8).Value = 2933.2222
$ws.Cells.Item(102, 9).Value = 2933.2222
$ws.Cells.Item(102, 11).Value = 2933.2222
$ws.Cells.Item(102, 13).Value = -1311.2222

$ws.Cells.Item(122, 8).Value = 1655.1
$ws.Cells.Item(122, 9).Value = 1403.4231
$ws.Cells.Item(122, 10).Value = 2122.5
$ws.Cells.Item(122, 11).Value = 4210.2693
$ws.Cells.Item(122, 12).Value = 6367.5
$ws.Cells.Item(122, 13).Value = -1760.2693
$ws.Cells.Item(122, 14).Value = -11267.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 71.46154
$ws.Cells.Item(7, 9).Value = 60.77778
$ws.Cells.Item(7, 10).Value = 95.5
$ws.Cells.Item(7, 11).Value = 60.77778
$ws.Cells.Item(7, 12).Value = 95.5
$ws.Cells.Item(7, 13).Value = 52.22222
$ws.Cells.Item(7, 14).Value = -321.5

$ws.Cells.Item(48, 8).Value = 5958.1665
$ws.Cells.Item(48, 10).Value = 5958.1665
$ws.Cells.Item(48, 12).Value = 5958.1665
$ws.Cells.Item(48, 14).Value = -6910.1665

$ws.Cells.Item(86, 8).Value = 265589.6
$ws.Cells.Item(86, 9).Value = 557578.25
$ws.Cells.Item(86, 10).Value = 2799.8
$ws.Cells.Item(86, 11).Value = 557578.25
$ws.Cells.Item(86, 12).Value = 2799.8
$ws.Cells.Item(86, 13).Value = -556455.25
$ws.Cells.Item(86, 14).Value = -5045.8

$ws.Cells.Item(89, 8).Value = 265589.6
$ws.Cells.Item(89, 9).Value = 557578.25
$ws.Cells.Item(89, 10).Value = 2799.8
$ws.Cells.Item(89, 11).Value = 2787891.25
$ws.Cells.Item(89, 12).Value = 13999
$ws.Cells.Item(89, 13).Value = -2782275.25
$ws.Cells.Item(89, 14).Value = -25231

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 64.2
$ws.Cells.Item(23, 9).Value = 23.333334
$ws.Cells.Item(23, 10).Value = 81.71429000000001
$ws.Cells.Item(23, 11).Value = 70.00000199999999
$ws.Cells.Item(23, 12).Value = 245.14287
$ws.Cells.Item(23, 13).Value = 164.999998
$ws.Cells.Item(23, 14).Value = -715.14287

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1970.1428
$ws.Cells.Item(113, 9).Value = 1927.75
$ws.Cells.Item(113, 10).Value = 2026.6666
$ws.Cells.Item(113, 11).Value = 1927.75
$ws.Cells.Item(113, 12).Value = 2026.6666
$ws.Cells.Item(113, 13).Value = 242.25
$ws.Cells.Item(113, 14).Value = -6366.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1315.5883
$ws.Cells.Item(46, 9).Value = 875.25
$ws.Cells.Item(46, 10).Value = 1451.0769
$ws.Cells.Item(46, 11).Value = 875.25
$ws.Cells.Item(46, 12).Value = 1451.0769
$ws.Cells.Item(46, 13).Value = -687.25
$ws.Cells.Item(46, 14).Value = -1827.0769

$ws.Cells.Item(132, 8).Value = 6971421.5
$ws.Cells.Item(132, 9).Value = 20904334
$ws.Cells.Item(132, 10).Value = 4964.6665
$ws.Cells.Item(132, 11).Value = 62713002
$ws.Cells.Item(132, 12).Value = 14893.9995
$ws.Cells.Item(132, 13).Value = -62710472
$ws.Cells.Item(132, 14).Value = -19953.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(51, 8).Value = 16696.666
$ws.Cells.Item(51, 9).Value = 15000
$ws.Cells.Item(51, 10).Value = 17545
$ws.Cells.Item(51, 11).Value = 15000
$ws.Cells.Item(51, 12).Value = 17545
$ws.Cells.Item(51, 13).Value = -14490
$ws.Cells.Item(51, 14).Value = -18565

$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(52, 9).Value = 0
$ws.Cells.Item(52, 11).Value = 0
$ws.Cells.Item(52, 13).ClearContents()

$ws.Cells.Item(61, 8).Value = 12525.5
$ws.Cells.Item(61, 9).Value = 10051
$ws.Cells.Item(61, 10).Value = 15000
$ws.Cells.Item(61, 11).Value = 10051
$ws.Cells.Item(61, 12).Value = 15000
$ws.Cells.Item(61, 13).Value = -9759
$ws.Cells.Item(61, 14).Value = -15584

$ws.Cells.Item(113, 8).Value = 33334760
$ws.Cells.Item(113, 9).Value = 71430360
$ws.Cells.Item(113, 10).Value = 1110
$ws.Cells.Item(113, 11).Value = 214291080
$ws.Cells.Item(113, 12).Value = 3330
$ws.Cells.Item(113, 13).Value = -214288910
$ws.Cells.Item(113, 14).Value = -7670

$ws.Cells.Item(126, 8).Value = 1706.7142
$ws.Cells.Item(126, 9).Value = 2000.4445
$ws.Cells.Item(126, 10).Value = 1178
$ws.Cells.Item(126, 11).Value = 6001.333500000001
$ws.Cells.Item(126, 12).Value = 3534
$ws.Cells.Item(126, 13).Value = -3531.333500000001
$ws.Cells.Item(126, 14).Value = -8474

$ws.Cells.Item(132, 8).Value = 2291.4138
$ws.Cells.Item(132, 9).Value = 1806.6522
$ws.Cells.Item(132, 10).Value = 4149.6665
$ws.Cells.Item(132, 11).Value = 5419.9566
$ws.Cells.Item(132, 12).Value = 12448.9995
$ws.Cells.Item(132, 13).Value = -2889.9566
$ws.Cells.Item(132, 14).Value = -17508.9995
